$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.029291292965256
$ws.Cells.Item(2, 4).Value = 1.037799921630124
$ws.Cells.Item(2, 5).Value = 1.029145161006346
$ws.Cells.Item(2, 6).Value = 1.046644144546457
$ws.Cells.Item(2, 9).Value = 1.035904970832552
$ws.Cells.Item(2, 10).Value = 1.034438881527789
$ws.Cells.Item(2, 11).Value = 1.040589773744273
$ws.Cells.Item(2, 12).Value = 1.03195992276475
$ws.Cells.Item(2, 13).Value = 1.049409003504707
$ws.Cells.Item(2, 14).Value = 1.005712725503983
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030242965046321
$ws.Cells.Item(3, 4).Value = 1.038519694917004
$ws.Cells.Item(3, 5).Value = 1.029953080002305
$ws.Cells.Item(3, 6).Value = 1.047504342387772
$ws.Cells.Item(3, 9).Value = 1.036093844866152
$ws.Cells.Item(3, 10).Value = 1.035031512372653
$ws.Cells.Item(3, 11).Value = 1.041119667576619
$ws.Cells.Item(3, 12).Value = 1.032575926150954
$ws.Cells.Item(3, 13).Value = 1.050080759443252
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030859121072944
$ws.Cells.Item(4, 4).Value = 1.038985409171182
$ws.Cells.Item(4, 5).Value = 1.030476558593249
$ws.Cells.Item(4, 6).Value = 1.04806123570434
$ws.Cells.Item(4, 9).Value = 1.036214393903319
$ws.Cells.Item(4, 10).Value = 1.035414743641325
$ws.Cells.Item(4, 11).Value = 1.041461847445733
$ws.Cells.Item(4, 12).Value = 1.032974572230354
$ws.Cells.Item(4, 13).Value = 1.050515076613365
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031118238422998
$ws.Cells.Item(5, 4).Value = 1.039181187475559
$ws.Cells.Item(5, 5).Value = 1.03069679544592
$ws.Cells.Item(5, 6).Value = 1.048295421041034
$ws.Cells.Item(5, 9).Value = 1.036264673372987
$ws.Cells.Item(5, 10).Value = 1.035575795480897
$ws.Cells.Item(5, 11).Value = 1.041605531911813
$ws.Cells.Item(5, 12).Value = 1.033142174117875
$ws.Cells.Item(5, 13).Value = 1.050697577546321
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03116175034369
$ws.Cells.Item(6, 4).Value = 1.039214059004586
$ws.Cells.Item(6, 5).Value = 1.030733783903807
$ws.Cells.Item(6, 6).Value = 1.048334745697068
$ws.Cells.Item(6, 9).Value = 1.036273092078408
$ws.Cells.Item(6, 10).Value = 1.035602833328644
$ws.Cells.Item(6, 11).Value = 1.041629647284611
$ws.Cells.Item(6, 12).Value = 1.033170315830864
$ws.Cells.Item(6, 13).Value = 1.050728215174883
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030862583077466
$ws.Cells.Item(7, 4).Value = 1.038988025203619
$ws.Cells.Item(7, 5).Value = 1.030479500756524
$ws.Cells.Item(7, 6).Value = 1.048064364637397
$ws.Cells.Item(7, 9).Value = 1.036215067310116
$ws.Cells.Item(7, 10).Value = 1.035416895855358
$ws.Cells.Item(7, 11).Value = 1.041463768025947
$ws.Cells.Item(7, 12).Value = 1.032976811692795
$ws.Cells.Item(7, 13).Value = 1.050517515539342
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.029612839767908
$ws.Cells.Item(8, 4).Value = 1.038043177137642
$ws.Cells.Item(8, 5).Value = 1.029418055053767
$ws.Cells.Item(8, 6).Value = 1.046934791937154
$ws.Cells.Item(8, 9).Value = 1.035969146033461
$ws.Cells.Item(8, 10).Value = 1.034639213093012
$ws.Cells.Item(8, 11).Value = 1.040768997467197
$ws.Cells.Item(8, 12).Value = 1.032168092858606
$ws.Cells.Item(8, 13).Value = 1.049636099057112
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.027413436026834
$ws.Cells.Item(9, 4).Value = 1.036378089972078
$ws.Cells.Item(9, 5).Value = 1.027553080729906
$ws.Cells.Item(9, 6).Value = 1.044946611163418
$ws.Cells.Item(9, 9).Value = 1.035523080256346
$ws.Cells.Item(9, 10).Value = 1.033267038848751
$ws.Cells.Item(9, 11).Value = 1.039539434864779
$ws.Cells.Item(9, 12).Value = 1.030743461919756
$ws.Cells.Item(9, 13).Value = 1.048080277515426
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025949107187538
$ws.Cells.Item(10, 4).Value = 1.035268020416693
$ws.Cells.Item(10, 5).Value = 1.026313494415641
$ws.Cells.Item(10, 6).Value = 1.043622772483552
$ws.Cells.Item(10, 9).Value = 1.035217192025961
$ws.Cells.Item(10, 10).Value = 1.032351105894797
$ws.Cells.Item(10, 11).Value = 1.038716236308036
$ws.Cells.Item(10, 12).Value = 1.029794062519784
$ws.Cells.Item(10, 13).Value = 1.047041355059323
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025315507986926
$ws.Cells.Item(11, 4).Value = 1.034787363058548
$ws.Cells.Item(11, 5).Value = 1.025777641146091
$ws.Cells.Item(11, 6).Value = 1.04304993713103
$ws.Cells.Item(11, 9).Value = 1.035082728947823
$ws.Cells.Item(11, 10).Value = 1.031954236082492
$ws.Cells.Item(11, 11).Value = 1.038358967889967
$ws.Cells.Item(11, 12).Value = 1.029383059025344
$ws.Cells.Item(11, 13).Value = 1.046591100658197
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025080231705015
$ws.Cells.Item(12, 4).Value = 1.034608828424744
$ws.Cells.Item(12, 5).Value = 1.025578737278106
$ws.Cells.Item(12, 6).Value = 1.042837221437511
$ws.Cells.Item(12, 9).Value = 1.035032481826055
$ws.Cells.Item(12, 10).Value = 1.031806782113394
$ws.Cells.Item(12, 11).Value = 1.038226140638361
$ws.Cells.Item(12, 12).Value = 1.029230408958728
$ws.Cells.Item(12, 13).Value = 1.046423798131006
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025130696054678
$ws.Cells.Item(13, 4).Value = 1.034647124555699
$ws.Cells.Item(13, 5).Value = 1.025621396662376
$ws.Cells.Item(13, 6).Value = 1.042882846898679
$ws.Cells.Item(13, 9).Value = 1.035043273639479
$ws.Cells.Item(13, 10).Value = 1.031838413233403
$ws.Cells.Item(13, 11).Value = 1.03825463803215
$ws.Cells.Item(13, 12).Value = 1.02926315222446
$ws.Cells.Item(13, 13).Value = 1.046459687695901
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025296058529548
$ws.Cells.Item(14, 4).Value = 1.034772605262605
$ws.Cells.Item(14, 5).Value = 1.025761196910507
$ws.Cells.Item(14, 6).Value = 1.043032352735582
$ws.Cells.Item(14, 9).Value = 1.035078581652787
$ws.Cells.Item(14, 10).Value = 1.031942048278552
$ws.Cells.Item(14, 11).Value = 1.038347990827674
$ws.Cells.Item(14, 12).Value = 1.029370440612838
$ws.Cells.Item(14, 13).Value = 1.046577272555164
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025397953197691
$ws.Cells.Item(15, 4).Value = 1.034849918494801
$ws.Cells.Item(15, 5).Value = 1.02584735050405
$ws.Cells.Item(15, 6).Value = 1.043124476323743
$ws.Cells.Item(15, 9).Value = 1.035100296146547
$ws.Cells.Item(15, 10).Value = 1.032005896130103
$ws.Cells.Item(15, 11).Value = 1.038405492450363
$ws.Cells.Item(15, 12).Value = 1.029436546538351
$ws.Cells.Item(15, 13).Value = 1.04664971281919
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025991166944755
$ws.Cells.Item(16, 4).Value = 1.035299920394223
$ws.Cells.Item(16, 5).Value = 1.0263490762177
$ws.Cells.Item(16, 6).Value = 1.043660798135525
$ws.Cells.Item(16, 9).Value = 1.035226073574165
$ws.Cells.Item(16, 10).Value = 1.032377439337842
$ws.Cells.Item(16, 11).Value = 1.03873992989088
$ws.Cells.Item(16, 12).Value = 1.029821341507023
$ws.Cells.Item(16, 13).Value = 1.047071228761348
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.02636339929746
$ws.Cells.Item(17, 4).Value = 1.035582198383338
$ws.Cells.Item(17, 5).Value = 1.026664036257729
$ws.Cells.Item(17, 6).Value = 1.043997325549994
$ws.Cells.Item(17, 9).Value = 1.035304432367971
$ws.Cells.Item(17, 10).Value = 1.032610428341025
$ws.Cells.Item(17, 11).Value = 1.038949495394438
$ws.Cells.Item(17, 12).Value = 1.030062738754388
$ws.Cells.Item(17, 13).Value = 1.047335529918548
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.026580560926408
$ws.Cells.Item(18, 4).Value = 1.035746847238372
$ws.Cells.Item(18, 5).Value = 1.026847833462661
$ws.Cells.Item(18, 6).Value = 1.044193654364167
$ws.Cells.Item(18, 9).Value = 1.035349943606837
$ws.Cells.Item(18, 10).Value = 1.032746301273009
$ws.Cells.Item(18, 11).Value = 1.039071652411173
$ws.Cells.Item(18, 12).Value = 1.030203550562808
$ws.Cells.Item(18, 13).Value = 1.047489654203462
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.026654615002632
$ws.Cells.Item(19, 4).Value = 1.03580298834973
$ws.Cells.Item(19, 5).Value = 1.026910518164935
$ws.Cells.Item(19, 6).Value = 1.044260603829215
$ws.Cells.Item(19, 9).Value = 1.035365428817244
$ws.Cells.Item(19, 10).Value = 1.032792626058943
$ws.Cells.Item(19, 11).Value = 1.039113291366192
$ws.Cells.Item(19, 12).Value = 1.030251565225936
$ws.Cells.Item(19, 13).Value = 1.047542200100742
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026323457649319
$ws.Cells.Item(20, 4).Value = 1.035551912525446
$ws.Cells.Item(20, 5).Value = 1.026630235085878
$ws.Cells.Item(20, 6).Value = 1.043961215395419
$ws.Cells.Item(20, 9).Value = 1.035296045284973
$ws.Cells.Item(20, 10).Value = 1.032585433468334
$ws.Cells.Item(20, 11).Value = 1.038927019162561
$ws.Cells.Item(20, 12).Value = 1.030036838188725
$ws.Cells.Item(20, 13).Value = 1.047307176838604
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025247361444647
$ws.Cells.Item(21, 4).Value = 1.034735654230428
$ws.Cells.Item(21, 5).Value = 1.025720025459954
$ws.Cells.Item(21, 6).Value = 1.042988325301336
$ws.Cells.Item(21, 9).Value = 1.03506819264187
$ws.Cells.Item(21, 10).Value = 1.031911531405764
$ws.Cells.Item(21, 11).Value = 1.03832050411155
$ws.Cells.Item(21, 12).Value = 1.029338846436341
$ws.Cells.Item(21, 13).Value = 1.046542648336026
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024571186352319
$ws.Cells.Item(22, 4).Value = 1.03422245776822
$ws.Cells.Item(22, 5).Value = 1.025148527438045
$ws.Cells.Item(22, 6).Value = 1.042376983522391
$ws.Cells.Item(22, 9).Value = 1.034923188300981
$ws.Cells.Item(22, 10).Value = 1.03148759811725
$ws.Cells.Item(22, 11).Value = 1.03793845953209
$ws.Cells.Item(22, 12).Value = 1.028900078200185
$ws.Cells.Item(22, 13).Value = 1.046061623981016
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024929600486484
$ws.Cells.Item(23, 4).Value = 1.034494510797103
$ws.Cells.Item(23, 5).Value = 1.025451414325973
$ws.Cells.Item(23, 6).Value = 1.042701033442464
$ws.Cells.Item(23, 9).Value = 1.035000222972269
$ws.Cells.Item(23, 10).Value = 1.031712354126014
$ws.Cells.Item(23, 11).Value = 1.038141055022722
$ws.Cells.Item(23, 12).Value = 1.029132668974384
$ws.Cells.Item(23, 13).Value = 1.046316655435953
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026341505422461
$ws.Cells.Item(24, 4).Value = 1.035565597398401
$ws.Cells.Item(24, 5).Value = 1.026645508113827
$ws.Cells.Item(24, 6).Value = 1.043977531902436
$ws.Cells.Item(24, 9).Value = 1.035299835646607
$ws.Cells.Item(24, 10).Value = 1.032596727654617
$ws.Cells.Item(24, 11).Value = 1.038937175448437
$ws.Cells.Item(24, 12).Value = 1.030048541512124
$ws.Cells.Item(24, 13).Value = 1.047319988492837
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027981696671209
$ws.Cells.Item(25, 4).Value = 1.036808563171526
$ws.Cells.Item(25, 5).Value = 1.028034569721284
$ws.Cells.Item(25, 6).Value = 1.045460325329268
$ws.Cells.Item(25, 9).Value = 1.035639902137671
$ws.Cells.Item(25, 10).Value = 1.033621985784326
$ws.Cells.Item(25, 11).Value = 1.039857926197521
$ws.Cells.Item(25, 12).Value = 1.031111705210238
$ws.Cells.Item(25, 13).Value = 1.048482801266039
